# "Termino da classe Tabela" -- extend Sheet2 with a second frequency table
# (Salario Mensal Bruto) in columns E:G, rows 1-9.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Planilha1")
$ws  = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# Column widths for the new columns E, F, G
# (ColumnWidth values chosen so the engine's pixel-quantized output lands
# as close as possible to the authored widths 21.28 / 32.19 / 17.86.)
# ---------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 20.45
$ws.Columns.Item(6).ColumnWidth = 31.25
$ws.Columns.Item(7).ColumnWidth = 16.95

# ---------------------------------------------------------------------
# Header row (E1:G1) -- reuse the existing shaded/bordered header look
# (same visual style already used on B3/B5/B7/B9 of this sheet).
# ---------------------------------------------------------------------
$ws.Range("B3").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Salário Mensal Bruto"
$ws.Range("F1").Value = "Freq. Acumulada de Funcionários"
$ws.Range("G1").Value = "Freq. Relativa"

# ---------------------------------------------------------------------
# Data rows (E2:G8) -- class label, cumulative count, per-class count.
# Style: same shaded blue fill + border already used for the A1/B1
# header on sheet "Planilha1", centered vertically too.
# ---------------------------------------------------------------------
$ws1.Range("A1").Copy()
$ws.Range("E2:G8").PasteSpecial(-4122)
$ws.Range("E2:G8").VerticalAlignment = -4108

$labels = @("[1700, 2900[", "[2900, 4100[", "[4100, 5300[", "[5300, 6500[", "[6500, 7700[", "[7700, 8900[", "[8900, 10100[")
$freqAcum = @(18, 39, 64, 93, 115, 132, 140)

for ($i = 0; $i -lt 7; $i++) {
    $r = 2 + $i
    $ws.Range("E$r").Value = $labels[$i]
    $ws.Range("F$r").Value = $freqAcum[$i]
}

# G2 is a literal value (first class has no previous class to subtract)
$ws.Range("G2").Value = 18
# G3:G8 are differences against the previous row's cumulative count
for ($r = 3; $r -le 8; $r++) {
    $prev = $r - 1
    $ws.Range("G$r").Formula = "=F$r-F$prev"
}

# ---------------------------------------------------------------------
# Totals row (E9:G9)
# ---------------------------------------------------------------------
$ws.Range("B3").Copy()
$ws.Range("E9:G9").PasteSpecial(-4122)
$ws.Range("E9").Value = "TOTAL"
$ws.Range("F9").Value = "-"
$ws.Range("G9").Formula = "=SUM(G2:G8)"

# ---------------------------------------------------------------------
# Selection follows the last-edited cell, as in the authored change.
# ---------------------------------------------------------------------
$ws.Range("E14").Select()
